$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3").Value = "MUZZIN K"
$ws.Rows("4:5").Delete()
